# "metadata for single view" -- add three new translation key/value rows
# (description_comment, project_affiliation, citation_format) to the
# translations sheet, right after the existing "Reuse?" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "description_comment"
$ws.Range("B39").Value = "Beschreibung / Kommentar"

$ws.Range("A40").Value = "project_affiliation"
$ws.Range("B40").Value = "Projektzugehörigkeit"

$ws.Range("A41").Value = "citation_format"
$ws.Range("B41").Value = "empfohlene Zitation"

# Re-establish the frozen header (rows 1-2) so the view keeps scrolling
# correctly now that the sheet has grown, then land the cursor just past
# the newly-entered rows (mirrors what a user ends up with after typing
# the three rows and pressing Enter).
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("A42").Select()
